$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.176.84'
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").Value = '1.629.55'
$ws.Range("E3").Value = '  -1.14%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = "'215.95"
$ws.Range("E5").Value = '  -0.68%  '
$ws.Range("E6").Value = '  +1.28%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").Value = "'0.255"
$ws.Range("E8").Value = '  -0.51%  '
$ws.Range("E9").Value = '  -1.04%  '
$ws.Range("D10").Value = "'20.27"
$ws.Range("E10").Value = '  +1.05%  '
$ws.Range("E11").Value = '  +0.58%  '
$ws.Range("D12").Value = '1.640.54'
$ws.Range("E12").Value = '  -0.57%  '
$ws.Range("E13").Value = '  -0.17%  '
$ws.Range("D14").Value = "'0.542"
$ws.Range("E14").Value = '  +0.16%  '
$ws.Range("D15").Value = '27.172.90'
$ws.Range("D16").Value = "'64.49"
$ws.Range("E16").Value = '  -4.77%  '
$ws.Range("D17").Value = '0.0₃0732'
$ws.Range("E17").Value = '  -1.19%  '
$ws.Range("D18").Value = "'215.79"
$ws.Range("E18").Value = '  -1.52%  '
$ws.Range("D20").Value = "'6.89"
$ws.Range("E20").Value = '  +0.34%  '
$ws.Range("E21").Value = '  -1.48%  '
$ws.Range("E22").Value = '  -0.85%  '
$ws.Range("E23").Value = '  -1.17%  '
$ws.Range("D24").Value = "'147.83"
$ws.Range("E24").Value = '  +0.23%  '
$ws.Range("E25").Value = '  +0.17%  '
$ws.Range("E26").Value = '  -3.79%  '
$ws.Range("E27").Value = '  -0.37%  '
$ws.Range("D28").Value = "'15.58"
$ws.Range("E28").Value = '  -1.50%  '
$ws.Range("E29").Value = '  -0.65%  '
$ws.Range("E30").Value = '  -0.33%  '
$ws.Range("D31").Value = "'3.39"
$ws.Range("E31").Value = '  +0.17%  '
$ws.Range("E32").Value = '  -1.26%  '
$ws.Range("D33").Value = '1.314.47'
$ws.Range("E33").Value = '  +3.88%  '
$ws.Range("E34").Value = '  -1.93%  '
$ws.Range("D35").Value = "'2.45"
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("E36").Value = '  -1.62%  '
$ws.Range("D37").Value = "'0.851"
$ws.Range("E37").Value = '  +1.13%  '
$ws.Range("E38").Value = '  -0.57%  '
$ws.Range("E39").Value = '  +0.06%  '
$ws.Range("E40").Value = '  +1.55%  '
$ws.Range("E41").Value = '  -0.95%  '
$ws.Range("D42").Value = "'63.64"
$ws.Range("E42").Value = '  +1.64%  '
$ws.Range("D43").Value = '1.767.68'
$ws.Range("E43").Value = '  -1.16%  '
$ws.Range("E44").Value = '  -4.31%  '
$ws.Range("D45").Value = "'90.69"
$ws.Range("E45").Value = '  -1.41%  '
$ws.Range("E46").Value = '  -0.68%  '
$ws.Range("D47").Value = '0.0₆0107'
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("D48").Value = "'0.804"
$ws.Range("E48").Value = '  +20.24%  '
$ws.Range("E49").Value = '  +0.68%  '
$ws.Range("D50").Value = "'7.56"
$ws.Range("E50").Value = '  -1.37%  '
$ws.Range("D51").Value = "'0.0953"
$ws.Range("E51").Value = '  -2.32%  '
